# The commit swaps the contents of ppt/theme/theme1.xml (used by the slide
# master -> "Integral" / "Red Violet" colours) and ppt/theme/theme2.xml
# (used by the notes master -> default "Office Theme" colours): after the
# edit, theme1.xml holds the default Office colour scheme and theme2.xml
# holds the old Integral/Red Violet colour scheme. Font scheme and format
# scheme (gradients/lines/effects) are already byte-identical between the
# two theme parts, so only the 12 theme colours actually need to change.
#
# Re-create that effect through the exposed PowerPoint object model by
# pushing the "Office Theme" RGB values onto the (single) slide master's
# theme colour scheme, in clrScheme slot order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
#   9 accent5, 10 accent6, 11 hlink, 12 folHlink

function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target ("Office Theme") colours, in a:clrScheme slot order.
$officeColors = @(
    @(0x00, 0x00, 0x00),  # 1  dk1
    @(0xFF, 0xFF, 0xFF),  # 2  lt1
    @(0x44, 0x54, 0x6A),  # 3  dk2
    @(0xE7, 0xE6, 0xE6),  # 4  lt2
    @(0x5B, 0x9B, 0xD5),  # 5  accent1
    @(0xED, 0x7D, 0x31),  # 6  accent2
    @(0xA5, 0xA5, 0xA5),  # 7  accent3
    @(0xFF, 0xC0, 0x00),  # 8  accent4
    @(0x44, 0x72, 0xC4),  # 9  accent5
    @(0x70, 0xAD, 0x47),  # 10 accent6
    @(0x05, 0x63, 0xC1),  # 11 hlink
    @(0x95, 0x4F, 0x72)   # 12 folHlink
)

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $rgb = $officeColors[$i]
    $colorScheme.Item($i + 1).RGB = RGBVal $rgb[0] $rgb[1] $rgb[2]
}
